$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 1.0.0 -> 0.1.0
$meta.Range("B3").Value = "0.1.0"

# Status: active -> draft
$meta.Range("B6").Value = "draft"

# Date: 2025-11-28T01:24:36+00:00 -> 2025-12-26T14:13:58+00:00
$meta.Range("B8").Value = "2025-12-26T14:13:58+00:00"

# Description: (empty) -> new extension description
$descText = "Extension to link goal evaluation observations to the patient goals being evaluated. Enables tracking of goal progress and outcomes over time."
$meta.Range("B11").Value = $descText

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")

# Definition (column M) of the root Extension row (row 2) gets the same new description text
$elements.Range("M2").Value = $descText
